$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (D, E, F) before the existing "Terms Typically Offered" column,
# shifting it to column G.
$ws.Columns("D:F").Insert()

# Header row
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

# Row 2
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "F, W"

# Row 3
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "F, W, SP"

# Row 4
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "F, W, SP"

# Row 5
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "F, W, SP"

# Row 6
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "F, W, SP"

# Row 7
$ws.Range("C7").Value = "one of the HLTH/KINE 250, HLTH/KINE 255, or HLTH/KINE 260."
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "HLTH 101."
$ws.Range("G7").Value = "F, W, SP "

# Row 8
$ws.Range("C8").Value = "HLTH/KINE 250, HLTH/KINE 255, or HLTH/KINE 260."
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "KINE 319 and STAT 218."
$ws.Range("G8").Value = "F, W, SP "

# Row 9
$ws.Range("C9").Value = "HLTH/KINE 250, HLTH/KINE 255, or HLTH/KINE 260."
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "FSN 210 and HLTH/KINE 265."
$ws.Range("G9").Value = "F, W, SP "

# Row 10
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "W, SP"

# Row 11
$ws.Range("C11").Value = "HLTH/KINE 265; and STAT 217 or STAT 218."
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "MCRO 221 or MCRO 224."
$ws.Range("G11").Value = "F, W "

# Row 12
$ws.Range("C12").Value = "HLTH 101; HLTH/KINE 298; and PSY 201 or PSY 202."
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "F, W"

# Row 13
$ws.Range("C13").Value = "BIO 231; and one of the HLTH/KINE 250; HLTH/KINE 255; HLTH/KINE 260; or HLTH/KINE 443."
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "BIO 232."
$ws.Range("G13").Value = "F, W "

# Row 14
$ws.Range("C14").Value = "BIO 231; BIO 232; and HLTH/KINE 265."
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "SP"

# Row 15
$ws.Range("C15").Value = "Completion of GE Area A with grades of C- or better; KINE 180 or HLTH/KINE 265; and one of the HLTH/KINE 250, HTLH/KINE 255, or HLTH/KINE 260."
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "HLTH 299."
$ws.Range("G15").Value = "W "

# Row 16
$ws.Range("C16").Value = "HLTH/KINE 298."
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "HLTH 299; and PSY 201 or PSY 202."
$ws.Range("G16").Value = "W, SP "

# Row 17
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "NA"
$ws.Range("G17").Value = "F, W, SP"

# Row 18
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = "F, W, SP"

# Row 19
$ws.Range("C19").Value = "HLTH 299; HLTH 334; and STAT 313."
$ws.Range("D19").Value = "NA"
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = "HLTH 310."
$ws.Range("G19").Value = "SP "

# Row 20
$ws.Range("C20").Value = "BIO 231; BIO 232; HLTH 299 or KINE 304; and HLTH 334 or KINE 266."
$ws.Range("D20").Value = "NA"
$ws.Range("E20").Value = "NA"
$ws.Range("F20").Value = "NA"
$ws.Range("G20").Value = "SP"

# Row 21
$ws.Range("C21").Value = "Junior standing; Completion of GE Areas D1, D2, D3, and D4/E; and ANT 360, or HLTH 298 and HLTH 334, or MCRO 221, or MCRO 224."
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "NA"
$ws.Range("G21").Value = "W"

# Row 22
$ws.Range("C22").Value = "HLTH/KINE 265; HLTH/KINE 298; and KINE 266 or HLTH 334."
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("F22").Value = "NA"
$ws.Range("G22").Value = "F, W, SP"

# Row 23
$ws.Range("C23").Value = "HLTH/KINE 320; and HLTH/KINE 434."
$ws.Range("D23").Value = "NA"
$ws.Range("E23").Value = "NA"
$ws.Range("F23").Value = "NA"
$ws.Range("G23").Value = "SP"

# Row 24
$ws.Range("D24").Value = "NA"
$ws.Range("E24").Value = "NA"
$ws.Range("F24").Value = "NA"
$ws.Range("G24").Value = "TBD"

# Row 25
$ws.Range("D25").Value = "NA"
$ws.Range("E25").Value = "NA"
$ws.Range("F25").Value = "NA"
$ws.Range("G25").Value = "F, W, SP"

# Row 26
$ws.Range("C26").Value = "one of the HLTH/KINE 250; HLTH/KINE 255; or HLTH/KINE 260; and KINE 266 or HLTH 334."
$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = "NA"
$ws.Range("F26").Value = "HLTH/KINE 320."
$ws.Range("G26").Value = "SP "

# Row 27
$ws.Range("C27").Value = "FSN 210 or KINE 451; and one of the HLTH/KINE 298, KINE 304, or FSN 310."
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = "NA"
$ws.Range("F27").Value = "HLTH 405; and one of the KINE 266, HLTH 334, or FSN 415."
$ws.Range("G27").Value = "F, W "

# Row 28
$ws.Range("C28").Value = "HLTH 402; KINE 320; KINE 434; completion of GE Area A with a grade of C- or better; completion of graduation writing requirement; and senior standing."
$ws.Range("D28").Value = "NA"
$ws.Range("E28").Value = "NA"
$ws.Range("F28").Value = "KINE 435."
$ws.Range("G28").Value = "F, W, SP "

# Row 29
$ws.Range("C29").Value = "HLTH 402; KINE 320; KINE 434; completion of GE Area A with a grade of C- or better; completion of graduation writing requirement; and senior standing."
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "KINE 435."
$ws.Range("G29").Value = "F, W, SP "

# Row 30
$ws.Range("C30").Value = "HLTH 402; KINE 320; KINE 434; completion of GE Area A with a grade of C- or better; completion of graduation writing requirement; senior standing and consent of instructor."
$ws.Range("D30").Value = "NA"
$ws.Range("E30").Value = "NA"
$ws.Range("F30").Value = "KINE 435."
$ws.Range("G30").Value = "F, W, SP "

# Row 31
$ws.Range("C31").Value = "Senior standing; completion of graduation writing requirement; minimum GPA of 3.0; KINE 434; and KINE 435."
$ws.Range("D31").Value = "NA"
$ws.Range("E31").Value = "NA"
$ws.Range("F31").Value = "NA"
$ws.Range("G31").Value = "F, W, SP"
